$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '42.593.26'
$ws.Range("E2").Value = '  +0.67%  '

# Row 3
$ws.Range("D3").Value = '2.287.57'
$ws.Range("E3").Value = '  -0.65%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.21%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.86'
$ws.Range("E5").Value = '  -0.64%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '105.17'
$ws.Range("E6").Value = '  +1.42%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.627'
$ws.Range("E7").Value = '  +0.43%  '

# Row 8
$ws.Range("E8").Value = '  -0.07%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.605'
$ws.Range("E9").Value = '  -0.56%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.62'
$ws.Range("E10").Value = '  -0.24%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0905'
$ws.Range("E11").Value = '  -0.12%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.41'
$ws.Range("E12").Value = '  +0.55%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.997'
$ws.Range("E14").Value = '  +3.57%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.27'
$ws.Range("E15").Value = '  -0.03%  '

# Row 16
$ws.Range("D16").Value = '2.634.31'
$ws.Range("E16").Value = '  -0.65%  '

# Row 17
$ws.Range("D17").Value = '2.278.76'
$ws.Range("E17").Value = '  -1.26%  '

# Row 18
$ws.Range("D18").Value = '42.586.04'
$ws.Range("E18").Value = '  +0.67%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.40'
$ws.Range("E19").Value = '  -0.82%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.88'
$ws.Range("E20").Value = '  +21.51%  '

# Row 21
$ws.Range("E21").Value = '  -0.88%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.81'
$ws.Range("E22").Value = '  +0.61%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.53'
$ws.Range("E23").Value = '  -0.15%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '265.36'
$ws.Range("E24").Value = '  -3.83%  '

# Row 25
$ws.Range("E25").Value = '  -2.94%  '

# Row 26
$ws.Range("E26").Value = '  +0.37%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.87'
$ws.Range("E27").Value = '  +0.11%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.28'
$ws.Range("E28").Value = '  +23.98%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.28'
$ws.Range("E29").Value = '  -2.93%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.47'
$ws.Range("E30").Value = '  -1.47%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '37.26'
$ws.Range("E31").Value = '  +1.69%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '166.91'
$ws.Range("E32").Value = '  +0.91%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0873'
$ws.Range("E33").Value = '  -0.14%  '

# Row 34
$ws.Range("E34").Value = '  -2.87%  '

# Row 35
$ws.Range("E35").Value = '  -1.01%  '

# Row 36
$ws.Range("E36").Value = '  -3.96%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.55'
$ws.Range("E37").Value = '  -0.48%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0352'
$ws.Range("E38").Value = '  -3.62%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.80'
$ws.Range("E39").Value = '  +2.31%  '

# Row 40
$ws.Range("E40").Value = '  -3.76%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.58'
$ws.Range("E41").Value = '  +4.86%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '70.54'

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.231'
$ws.Range("E43").Value = '  +2.08%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '94.54'
$ws.Range("E44").Value = '  -0.19%  '

# Row 45
$ws.Range("E45").Value = '  -0.11%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.20'
$ws.Range("E46").Value = '  +0.86%  '

# Row 47
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '1.731.84'
$ws.Range("E47").Value = '  +8.85%  '

# Row 48
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '113.63'
$ws.Range("E48").Value = '  +0.56%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '79.51'
$ws.Range("E49").Value = '  -2.43%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.74'
$ws.Range("E50").Value = '  -2.41%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.17'
$ws.Range("E51").Value = '  -0.68%  '
